$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 71.333336
$ws.Range("I9").Value = 74
$ws.Range("J9").Value = 70
$ws.Range("K9").Value = 74
$ws.Range("L9").Value = 70
$ws.Range("M9").Value = 95
$ws.Range("N9").Value = -408
$ws.Range("H40").Value = 2500
$ws.Range("J40").Value = 2500
$ws.Range("L40").Value = 2500
$ws.Range("N40").Value = -2850
$ws.Range("H58").Value = 6455.6665
$ws.Range("I58").Value = 1400
$ws.Range("J58").Value = 8983.5
$ws.Range("K58").Value = 4200
$ws.Range("L58").Value = 26950.5
$ws.Range("M58").Value = -4050
$ws.Range("N58").Value = -27250.5
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").Value = $null
$ws.Range("H80").Value = 749.5
$ws.Range("H83").Value = 749.5
$ws.Range("H112").Value = 1648.7727
$ws.Range("J112").Value = 1788.0526
$ws.Range("L112").Value = 5364.1578
$ws.Range("N112").Value = -7580.1578
$ws.Range("H116").Value = 7311.857
$ws.Range("I116").Value = 7002
$ws.Range("K116").Value = 7002
$ws.Range("M116").Value = -3560
$ws.Range("H138").Value = 2866.0938
$ws.Range("J138").Value = 3010.434
$ws.Range("L138").Value = 9031.302
$ws.Range("N138").Value = -19311.302

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1514.9375
$ws.Range("I2").Value = 1555.9333
$ws.Range("K2").Value = 1555.9333
$ws.Range("M2").Value = -1442.9333
$ws.Range("H32").Value = 4762.1
$ws.Range("I32").Value = 3132.4856
$ws.Range("J32").Value = 16169.4
$ws.Range("K32").Value = 3132.4856
$ws.Range("L32").Value = 16169.4
$ws.Range("M32").Value = -2845.4856
$ws.Range("N32").Value = -16743.4
$ws.Range("H45").Value = 2511
$ws.Range("I45").Value = 2511
$ws.Range("K45").Value = 2511
$ws.Range("M45").Value = -2134
$ws.Range("H97").Value = 303.72726
$ws.Range("I97").Value = 234.1
$ws.Range("K97").Value = 234.1
$ws.Range("M97").Value = 261.9
$ws.Range("H102").Value = 1457.6666
$ws.Range("I102").Value = 1457.6666
$ws.Range("K102").Value = 1457.6666
$ws.Range("M102").Value = 164.3334
$ws.Range("H110").Value = 1200
$ws.Range("I110").Value = 1200
$ws.Range("K110").Value = 1200
$ws.Range("M110").Value = 845
$ws.Range("H116").Value = 1514.9375
$ws.Range("I116").Value = 1555.9333
$ws.Range("K116").Value = 1555.9333
$ws.Range("M116").Value = 738.0667000000001
$ws.Range("H122").Value = 5648.9
$ws.Range("I122").Value = 5943.222
$ws.Range("K122").Value = 17829.666
$ws.Range("M122").Value = -15379.666
$ws.Range("H132").Value = 3401.7144
$ws.Range("I132").Value = 2772.3076
$ws.Range("K132").Value = 8316.9228
$ws.Range("M132").Value = -5786.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1514.9375
$ws.Range("I3").Value = 1555.9333
$ws.Range("K3").Value = 1555.9333
$ws.Range("M3").Value = -1441.9333
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 4377.8
$ws.Range("I2").Value = 440
$ws.Range("J2").Value = 7003
$ws.Range("K2").Value = 440
$ws.Range("L2").Value = 7003
$ws.Range("M2").Value = -327
$ws.Range("N2").Value = -7229
$ws.Range("H8").Value = 7170
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 7170
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 7170
$ws.Range("M8").Value = $null
$ws.Range("N8").Value = -7450
$ws.Range("H16").Value = 2226.5
$ws.Range("I16").Value = 2226.5
$ws.Range("K16").Value = 2226.5
$ws.Range("M16").Value = -1939.5
$ws.Range("H86").Value = 6603.864
$ws.Range("I86").Value = 7997.364
$ws.Range("J86").Value = 5210.364
$ws.Range("K86").Value = 7997.364
$ws.Range("L86").Value = 5210.364
$ws.Range("M86").Value = -6874.364
$ws.Range("N86").Value = -7456.364
$ws.Range("H89").Value = 6603.864
$ws.Range("I89").Value = 7997.364
$ws.Range("J89").Value = 5210.364
$ws.Range("K89").Value = 39986.82
$ws.Range("L89").Value = 26051.82
$ws.Range("M89").Value = -34370.82
$ws.Range("N89").Value = -37283.82
$ws.Range("H107").Value = 2227.9167
$ws.Range("I107").Value = 1164.3334
$ws.Range("K107").Value = 1164.3334
$ws.Range("M107").Value = 755.6666
$ws.Range("H113").Value = 2226.5
$ws.Range("I113").Value = 2226.5
$ws.Range("K113").Value = 2226.5
$ws.Range("M113").Value = -56.5
$ws.Range("H132").Value = 2311.84
$ws.Range("I132").Value = 1711.0555
$ws.Range("K132").Value = 5133.166499999999
$ws.Range("M132").Value = -2603.166499999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null
$ws.Range("H107").Value = 152.875
$ws.Range("J107").Value = 160.42857
$ws.Range("L107").Value = 481.28571
$ws.Range("N107").Value = -4321.28571
$ws.Range("H113").Value = 2127.2
$ws.Range("I113").Value = 2375
$ws.Range("J113").Value = 2089.077
$ws.Range("K113").Value = 7125
$ws.Range("L113").Value = 6267.231000000001
$ws.Range("M113").Value = -4955
$ws.Range("N113").Value = -10607.231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 15065.5
$ws.Range("I31").Value = 131
$ws.Range("J31").Value = 30000
$ws.Range("K31").Value = 131
$ws.Range("L31").Value = 30000
$ws.Range("M31").Value = 161
$ws.Range("N31").Value = -30584
$ws.Range("H37").Value = 15065.5
$ws.Range("I37").Value = 131
$ws.Range("J37").Value = 30000
$ws.Range("K37").Value = 131
$ws.Range("L37").Value = 30000
$ws.Range("M37").Value = 146
$ws.Range("N37").Value = -30554
$ws.Range("H46").Value = 34942
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 34942
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 34942
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = -35254
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = $null
$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 5000
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4730
$ws.Range("N70").Value = $null
$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 5000
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -4064
$ws.Range("N73").Value = $null
$ws.Range("H122").Value = 1666.3334
$ws.Range("I122").Value = 1124.5
$ws.Range("K122").Value = 3373.5
$ws.Range("M122").Value = -923.5
$ws.Range("H123").Value = 20000
$ws.Range("J123").Value = 20000
$ws.Range("L123").Value = 20000
$ws.Range("N123").Value = -24900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8000
$ws.Range("J22").Value = 8000
$ws.Range("L22").Value = 8000
$ws.Range("N22").Value = -8590
$ws.Range("H27").Value = 8000
$ws.Range("J27").Value = 8000
$ws.Range("L27").Value = 8000
$ws.Range("N27").Value = -8214
$ws.Range("H46").Value = 1422.75
$ws.Range("I46").Value = 1422.75
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1422.75
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1234.75
$ws.Range("N46").Value = $null
$ws.Range("H55").Value = 229.33333
$ws.Range("I55").Value = 300.5
$ws.Range("K55").Value = 300.5
$ws.Range("M55").Value = -127.5
$ws.Range("H132").Value = 3942.4
$ws.Range("I132").Value = 3748
$ws.Range("K132").Value = 11244
$ws.Range("M132").Value = -8714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 16398.6
$ws.Range("I81").Value = 14666.333
$ws.Range("K81").Value = 29332.666
$ws.Range("M81").Value = -28271.666
$ws.Range("H84").Value = 16398.6
$ws.Range("I84").Value = 14666.333
$ws.Range("K84").Value = 146663.33
$ws.Range("M84").Value = -141359.33
